$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16
$ws.Cells.Item(16, 1).Value = "Invalid Date Error"
$ws.Cells.Item(16, 2).Value = "date must be a ``date`` type, but the final value was: ``Invalid Date``."

# Row 17
$ws.Cells.Item(17, 1).Value = "Invalid Visit Reason"
$ws.Cells.Item(17, 2).Value = "Field is required"

# Row 18
$ws.Cells.Item(18, 1).Value = "Invalid Birthdate"
$ws.Cells.Item(18, 2).Value = "birthdate must be a ``date`` type, but the final value was: ``Invalid Date``."

# Update selection to match new active cell
$ws.Range("D20").Select()
